# Applies the "review" fix-up described in the commit:
#  - Column C (difference) is recomputed as Congruent - Incongruent (A-B)
#    instead of the inconsistent B-A/A-B mix that was there before.
#  - A new column D computes the squared deviation of each difference
#    from the mean difference: (diff - diff_mean)^2
#  - The old, partly-broken "E" helper column (diff mean / SE / Tstat /
#    df / Tcritical) is removed and replaced with a new "G" helper
#    column that derives the t-statistic from first principles:
#       diff_mean -> sum of squared deviations -> S -> SE -> t_stat

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inferential Stat. Analysis")

# ---------------------------------------------------------------------
# 1. New header for column D
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "(diff-diff_mean)^2"

# Match the formatting already used for the other header cell (C1).
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "(diff-diff_mean)^2"

# ---------------------------------------------------------------------
# 2. Remove the old "E" helper column entirely
# ---------------------------------------------------------------------
$ws.Range("E3").Clear()
$ws.Range("E5").Clear()
$ws.Range("E6").Clear()
$ws.Range("E8").Clear()
$ws.Range("E9").Clear()
$ws.Range("E11").Clear()
$ws.Range("E12").Clear()
$ws.Range("E14").Clear()
$ws.Range("E15").Clear()
$ws.Range("E18").Clear()
$ws.Range("E19").Clear()
$ws.Range("E21").Clear()
$ws.Range("E22").Clear()

# D23 used to be an empty, styled placeholder cell - drop the formatting
# since the column is now fully populated with real formulas.
$ws.Range("D23").ClearFormats()

# ---------------------------------------------------------------------
# 3. Column C: congruent - incongruent, as one shared formula C2:C25
# ---------------------------------------------------------------------
$ws.Range("C2:C25").Formula = "=A2-B2"

# ---------------------------------------------------------------------
# 4. Column G: diff_mean, sum of squared deviations, S, SE, t_stat
# ---------------------------------------------------------------------
$ws.Range("E5").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "diff_mean"

$ws.Range("G4").Formula = "=AVERAGE(C2:C25)"

$ws.Range("E8").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = "sum of dev"

$ws.Range("G8").Formula = "=SUM(D2:D25)"
$ws.Range("E5").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Formula = "=SUM(D2:D25)"

$ws.Range("E11").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("G11").Value = "Squared dev diff  S"

$ws.Range("G12").Formula = "=SQRT(`$G`$8/24)"

$ws.Range("E14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value = "S/sqrt(n)"

$ws.Range("G15").Formula = "=`$G`$12/SQRT(24)"

$ws.Range("E18").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G17").Value = "t_stat"

$ws.Range("G18").Formula = "=(`$A`$30 - `$A`$33)/(`$G`$15)"

# ---------------------------------------------------------------------
# 5. Column D: squared deviation of each difference from the mean,
#    as a shared formula D2:D25 (depends on G4, so fill in after G4).
# ---------------------------------------------------------------------
$ws.Range("D2:D25").Formula = "=(C2-`$G`$4)^2"

# ---------------------------------------------------------------------
# 6. A handful of cells that only carry over formatting (no content)
#    in the reviewed version - reuse the existing blank-but-styled
#    cells as the format source.
# ---------------------------------------------------------------------
$ws.Range("D26").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").ClearContents()

$ws.Range("D26").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").ClearContents()

$ws.Range("D26").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("G23").ClearContents()

$ws.Range("D26").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G27").ClearContents()

$ws.Range("D26").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").ClearContents()

$ws.Range("D26").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").ClearContents()

$excel.CutCopyMode = 0
